$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K (strikeout) value for column G,
# regenerated from the updated save_data source (use K instead of Strike#).
$kValues = @{
    2 = 2
    3 = 2
    4 = 1
    5 = 2
    6 = 1
    7 = 3
    8 = 3
    9 = 0
    10 = 2
    11 = 2
    12 = 0
    14 = 1
    15 = 2
    16 = 2
    17 = 1
    18 = 1
    19 = 2
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 2
    25 = 2
    26 = 1
    27 = 0
    28 = 1
    29 = 0
    30 = 1
    31 = 1
    32 = 1
    33 = 2
    34 = 2
    35 = 1
    36 = 3
    37 = 1
    38 = 2
    39 = 1
    40 = 0
    41 = 1
    42 = 1
    43 = 2
    44 = 3
    45 = 2
    46 = 1
    47 = 0
    48 = 2
    49 = 2
    50 = 1
    51 = 0
    52 = 0
    53 = 3
    54 = 2
    55 = 1
    56 = 0
    57 = 2
    58 = 3
    59 = 0
    60 = 1
    61 = 1
    62 = 1
    63 = 1
    64 = 2
    65 = 2
    66 = 1
    67 = 1
    68 = 0
    69 = 3
    70 = 1
    71 = 1
    72 = 1
    73 = 1
    75 = 2
    76 = 0
    77 = 3
    78 = 1
    79 = 1
    80 = 2
    81 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

Write-Host "Updated" $kValues.Count "K values in column G"
